# Update TPM-derived NATMI ligand/receptor metrics on Sheet1 (Ceacam1-Sele)
# to reflect newly recomputed TPM values (commit: "update scripts wuth new tpm").
# Only numeric value cells in columns G:J, M:T (rows 2-9) change; identifiers,
# detection-rate counts (E:F, K:L) and headers are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 37.08402533333333
$ws.Range("H2").Value = 111.252076
$ws.Range("I2").Value = 0.8732570105511104
$ws.Range("J2").Value = 0.8732570105511105
$ws.Range("M2").Value = 12.67919733333333
$ws.Range("N2").Value = 38.037592
$ws.Range("O2").Value = 0.9871416146107245
$ws.Range("P2").Value = 0.9871416146107247
$ws.Range("Q2").Value = 470.1956751156658
$ws.Range("R2").Value = 4231.761076040992
$ws.Range("S2").Value = 0.8620283353655576
$ws.Range("T2").Value = 0.8620283353655579
$ws.Range("G3").Value = 37.08402533333333
$ws.Range("H3").Value = 111.252076
$ws.Range("I3").Value = 0.8732570105511104
$ws.Range("J3").Value = 0.8732570105511105
$ws.Range("M3").Value = 0.1651576666666667
$ws.Range("N3").Value = 0.495473
$ws.Range("O3").Value = 0.01285838538927542
$ws.Range("P3").Value = 0.01285838538927542
$ws.Range("Q3").Value = 6.124711094660888
$ws.Range("R3").Value = 55.122399851948
$ws.Range("S3").Value = 0.01122867518555273
$ws.Range("T3").Value = 0.01122867518555273
$ws.Range("I4").Value = 0.003499940937530227
$ws.Range("J4").Value = 0.003499940937530227
$ws.Range("M4").Value = 12.67919733333333
$ws.Range("N4").Value = 38.037592
$ws.Range("O4").Value = 0.9871416146107245
$ws.Range("P4").Value = 0.9871416146107247
$ws.Range("Q4").Value = 1.884504873254222
$ws.Range("R4").Value = 16.960543859288
$ws.Range("S4").Value = 0.003454937348115761
$ws.Range("T4").Value = 0.003454937348115762
$ws.Range("I5").Value = 0.003499940937530227
$ws.Range("J5").Value = 0.003499940937530227
$ws.Range("M5").Value = 0.1651576666666667
$ws.Range("N5").Value = 0.495473
$ws.Range("O5").Value = 0.01285838538927542
$ws.Range("P5").Value = 0.01285838538927542
$ws.Range("S5").Value = 0.00004500358941446558
$ws.Range("T5").Value = 0.00004500358941446558
$ws.Range("G6").Value = 0.5778596666666668
$ws.Range("H6").Value = 1.733579
$ws.Range("I6").Value = 0.01360747654807074
$ws.Range("J6").Value = 0.01360747654807074
$ws.Range("M6").Value = 12.67919733333333
$ws.Range("N6").Value = 38.037592
$ws.Range("O6").Value = 0.9871416146107245
$ws.Range("P6").Value = 0.9871416146107247
$ws.Range("Q6").Value = 7.326796744640891
$ws.Range("R6").Value = 65.94117070176802
$ws.Range("S6").Value = 0.01343250637044012
$ws.Range("T6").Value = 0.01343250637044013
$ws.Range("G7").Value = 0.5778596666666668
$ws.Range("H7").Value = 1.733579
$ws.Range("I7").Value = 0.01360747654807074
$ws.Range("J7").Value = 0.01360747654807074
$ws.Range("M7").Value = 0.1651576666666667
$ws.Range("N7").Value = 0.495473
$ws.Range("O7").Value = 0.01285838538927542
$ws.Range("P7").Value = 0.01285838538927542
$ws.Range("Q7").Value = 0.09543795420744447
$ws.Range("R7").Value = 0.8589415878670001
$ws.Range("S7").Value = 0.0001749701776306207
$ws.Range("T7").Value = 0.0001749701776306207
$ws.Range("G8").Value = 4.655821
$ws.Range("H8").Value = 13.967463
$ws.Range("I8").Value = 0.1096355719632885
$ws.Range("J8").Value = 0.1096355719632885
$ws.Range("M8").Value = 12.67919733333333
$ws.Range("N8").Value = 38.037592
$ws.Range("O8").Value = 0.9871416146107245
$ws.Range("P8").Value = 0.9871416146107247
$ws.Range("Q8").Value = 59.03207320767734
$ws.Range("R8").Value = 531.288658869096
$ws.Range("S8").Value = 0.1082258355266109
$ws.Range("T8").Value = 0.108225835526611
$ws.Range("G9").Value = 4.655821
$ws.Range("H9").Value = 13.967463
$ws.Range("I9").Value = 0.1096355719632885
$ws.Range("J9").Value = 0.1096355719632885
$ws.Range("M9").Value = 0.1651576666666667
$ws.Range("N9").Value = 0.495473
$ws.Range("O9").Value = 0.01285838538927542
$ws.Range("P9").Value = 0.01285838538927542
$ws.Range("Q9").Value = 0.7689445327776667
$ws.Range("R9").Value = 6.920500794999001
$ws.Range("S9").Value = 0.001409736436677603
$ws.Range("T9").Value = 0.001409736436677603
